$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(10)
$r = $p.Range
$r.Collapse(0)
for ($i = 0; $i -lt 4; $i++) {
    $r.Text = "`r"
    $r.Collapse(0)
}

$notePara = $d.Paragraphs.Item(12)
Write-Output "notePara: start=$($notePara.Range.Start) end=$($notePara.Range.End)"

$xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:spacing w:after="0" w:line="360" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Note: Feature from Part C: Press the C button to change from overhead view to 3</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:vertAlign w:val="superscript"/>
    </w:rPr>
    <w:t>rd</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> person 3D view</w:t>
  </w:r>
</w:p>
"@
$notePara.Range.InsertXML($xml)
Write-Output "done"
